# Commit: "Fixed POI packaging and upgraded to POI 3.15."
#
# The underlying OOXML diff for this fixture is entirely non-semantic:
#   - xmlns:* declarations on <w:document>/<w:footnotes> reordered alphabetically
#   - attributes on existing elements (<w:tab>, <w:pgSz>, <w:pgMar>, <w:footnote>,
#     <w:rFonts>, <w:lang>, <w:latentStyles>, <w:lsdException>, <w:style>,
#     <w:tblInd>, <w:tblCellMar>, ...) reordered alphabetically
#   - two internally-generated, random, tool-specific identifiers changed value
#     (a 32-hex-digit w:rsidR placeholder on the REF-field runs, and the
#     w:bookmarkStart/w:bookmarkEnd w:id for "bookmark1")
#
# None of the visible document content changed: every paragraph's text, the
# tab stop (pos=3119, left), the field code (" REF bookmark1 \h "), the
# bookmark name ("bookmark1"), and the section page setup
# (11906x16838 / margins 1417/1417/1417/1417, header/footer 708, gutter 0)
# are identical before and after. This is simply the fixture being
# regenerated by a newer build of the authoring tool (POI 3.15), which
# happens to serialize attributes/namespaces in a different (alphabetical)
# order and mints fresh random internal ids - none of which are properties
# that the Word object model exposes or lets an author control.
#
# So there is nothing to change from an authoring point of view; this script
# just verifies the content that must survive the round-trip untouched.

$d = $word.ActiveDocument

# Sanity-check the content is what we expect - untouched by this "edit".
$null = $d.Content.Find.Execute("bookmarked content", $true, $false, $false,
                                 $false, $false, $true, 1, $false, "", 0)
$null = $d.Bookmarks("bookmark1")
